# Apply JIRA id / description updates to the "Test Cases" sheet (sheet1)
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# Row 2 (ENWIAM001 registration test case): append new OPQA ids + new bullet to description
$ws.Range("B2").Value = "OPQA-1719||`nOPQA-1676||OPQA-1744||`nOPQA-1760||OPQA-1763`n||OPQA-1766||OPQA-2038||OPQA-2359||OPQA-2139"
$ws.Range("C2").Value = 'Verify that ENW registration screen should be displayed and User should be able to enter email address (required), name (required), and password (required).||Verify that "Sign up" link should be displayed on ENW registration page .||Verify that the user should be able click on "sign up" button after filling the above fields correctly.||Verify that user should get an Email verification Link on the registered Email Id .||Verify that after clicking verification link user should get the message as" Success!You have successfully activated your account. Please sign in."||Verify that after completion of verification process,user should be able to sign into ENW ||Verify that after successful registration on the ENW landing screen using Facebook, users who already has Steam account with the same email address are prompted to link their Steam account with the newly created Facebook account.'

# Row 4 (password strength test case): append new OPQA id + new bullet to description
$ws.Range("B4").Value = 'OPQA-1741||OPQA-2004'
$ws.Range("C4").Value = 'Verify that Passwords should be at least 8 characters,Must contain at least 1 lowercase letter,Must contain at least 1 uppercase letter,Must contain at least 1 number,should have at least one special character and must not be empty.||Verify that The system shall not allow a user to sign in to Neon with STeAM credentials that are locked.'

# Row 11 (STeAM sign-in test case): append new OPQA id + new bullet to description
$ws.Range("B11").Value = 'OPQA-2007||OPQA-3652||OPQA-2008||OPQA-2009||OPQA-3333'
$ws.Range("C11").Value = 'Verify that STeAM user is able to submit an email address and password on the ENW Landing screen.||Verify that,user should receive the ENW EULA acceptance after signed into ENW for the first time.||Verify that A user shall successfully authenticate by supplying correct STeAM credentials (email address + password), on the ENW landing screen.||Verify that A user should not be allowed to sign-in to ENW if an incorrect email address and password combination is provided on the ENW landing screen||Verify that in account setting page,after clicking on "Link account"Button, "Enter your existing account credentials (CortellisTM,EndNoteTM Online,InCitesTM,ResearcherID,Thomson InnovationTM,Web of ScienceTM)to link your accounts" text should be displayed'

# Row 18 (account lockout test case): append new OPQA ids + new bullets to description
$ws.Range("B18").Value = 'OPQA-2001||OPQA-2005||OPQA-1870'
$ws.Range("C18").Value = 'Verify that A user signing in to NEON or ENW using STeAM shall be locked out of their account after 10 invalid attempts.||Verify that The system shall not allow a user to sign in to ENW with STeAM credentials that are locked.||Verify that If the STeAM account that is trying to be linked/merged by the user is in a "locked" status, then the link/merge shall not be made and the user shall be informed that the STeAM account is locked.'

# Update the view/selection to match the saved workbook state (scrolled to row 4, C5 selected)
$ws.Activate()
$ws.Application.ActiveWindow.ScrollRow = 4
$ws.Range("C5").Select()

